$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 39 (shifts existing rows 39-66 down to 41-68)
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(39).Insert()

# Row 39 (new)
$ws.Range("A39").Value = 8
$ws.Range("B39").Value = "Terminal La Palmera de La Serena"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44897
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103003
$ws.Range("J39").Value = "Damasco"
$ws.Range("K39").Value = "Castle Brite"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 100
$ws.Range("N39").Value = 25000
$ws.Range("O39").Value = 26000
$ws.Range("P39").Value = 25500
$ws.Range("Q39").Value = "`$/caja 18 kilos"
$ws.Range("R39").Value = "Región de O'Higgins"
$ws.Range("S39").Value = 1417
$ws.Range("T39").Value = 18

# Row 40 (new)
$ws.Range("A40").Value = 8
$ws.Range("B40").Value = "Terminal La Palmera de La Serena"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44897
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100103
$ws.Range("H40").Value = "Frutos de hueso (carozo)"
$ws.Range("I40").Value = 100103003
$ws.Range("J40").Value = "Damasco"
$ws.Range("K40").Value = "Castle Brite"
$ws.Range("L40").Value = "Segunda"
$ws.Range("M40").Value = 160
$ws.Range("N40").Value = 23000
$ws.Range("O40").Value = 24000
$ws.Range("P40").Value = 23500
$ws.Range("Q40").Value = "`$/caja 18 kilos"
$ws.Range("R40").Value = "Región de O'Higgins"
$ws.Range("S40").Value = 1306
$ws.Range("T40").Value = 18
